# Removed Test Case Inter-Dependency
# - Update the product name text used on both sheets (new "-1st" suffix instead of "-FEE-FLAT")
# - Change the shortname value on ProductLoanInput from the numeric 2617 to the text "261z"
# - Make ProductLoanOutput the active/selected sheet (was ProductLoanInput), with B1 selected
# - Leave ProductLoanInput's selection on B1 as well (previously B15)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$newProductName = "2617-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-DISBURSE-1st"

$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

$ws1.Range("B2").Value = "261z"

$ws1.Range("B1").Select()
$ws2.Activate()
